# Atualizando o arquivo XLSX
# Apply updated odds/statistics values to Sheet1 as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "K2" = 13
    "G4" = 1.65
    "I4" = 4.75
    "AF4" = 26
    "G5" = 2.05
    "H5" = 4.1
    "I5" = 2.88
    "P5" = 1.22
    "Q5" = 4
    "AC5" = 29
    "G6" = 1.7
    "H6" = 4.33
    "I6" = 4
    "R6" = 1.62
    "S6" = 2.2
    "U6" = 9.5
    "X6" = 12
    "G10" = 1.7
    "H10" = 4
    "I10" = 4.5
    "X10" = 13
    "AA10" = 7.5
    "AF10" = 23
    "AI10" = 34
    "AJ10" = 34
    "G11" = 1.65
    "I11" = 5.25
    "N11" = 1.9
    "O11" = 1.9
    "J14" = 1.06
    "K14" = 10
    "AA14" = 7
    "AE14" = 13
    "G15" = 1.67
    "I15" = 4.33
    "U15" = 8.5
    "J16" = 1.03
    "K16" = 10
    "H17" = 4.5
    "J17" = 19
    "N17" = 1.5
    "O17" = 2.5
    "T17" = 10
    "V17" = 9
    "Z17" = 19
    "AA17" = 9.5
    "AH17" = 67
    "G18" = 2.05
    "I18" = 3.6
    "J18" = 1.06
    "K18" = 10
    "L18" = 1.3
    "M18" = 3.4
    "N18" = 2
    "O18" = 1.8
    "U18" = 9.5
    "AG18" = 13
    "J20" = 1.05
    "K20" = 11
    "N20" = 2
    "O20" = 1.8
    "J23" = 1.06
    "K23" = 10
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
